$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers ("1.002", "310.47", ...) need the
# NumberFormat forced to Text first, otherwise Excel will store them as numeric values
# instead of the text strings the source data actually contains.
foreach ($addr in @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D22", "D23", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D49", "D50", "D51")) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.209.21"
$ws.Range("E2").Value = "  -3.79%  "
$ws.Range("D3").Value = "1.809.39"
$ws.Range("E3").Value = "  -3.77%  "
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "310.47"
$ws.Range("E5").Value = "  -1.72%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "0.4208"
$ws.Range("E7").Value = "  -2.33%  "
$ws.Range("D8").Value = "0.3557"
$ws.Range("E8").Value = "  -4.61%  "
$ws.Range("D9").Value = "0.07113"
$ws.Range("E9").Value = "  -4.22%  "
$ws.Range("D10").Value = "0.8512"
$ws.Range("E10").Value = "  -3.90%  "
$ws.Range("D11").Value = "20.16"
$ws.Range("E11").Value = "  -4.75%  "
$ws.Range("D12").Value = "1.826.95"
$ws.Range("E12").Value = "  -6.79%  "
$ws.Range("D13").Value = "5.311"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").Value = "6.369"
$ws.Range("E14").Value = "  -4.00%  "
$ws.Range("D15").Value = "0.06840"
$ws.Range("E15").Value = "  -2.26%  "
$ws.Range("D16").Value = "1.004"
$ws.Range("E16").Value = "  +0.05%  "
$ws.Range("D17").Value = "80.94"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "0.000008772"
$ws.Range("E18").Value = "  -4.09%  "
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("D20").Value = "15.13"
$ws.Range("E20").Value = "  -3.35%  "
$ws.Range("D21").Value = "27.340.55"
$ws.Range("E21").Value = "  -3.87%  "
$ws.Range("D22").Value = "5.112"
$ws.Range("E22").Value = "  +0.44%  "
$ws.Range("D23").Value = "10.90"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").Value = "2.057.82"
$ws.Range("E24").Value = "  -5.95%  "
$ws.Range("D25").Value = "1.970"
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("D26").Value = "153.59"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("D27").Value = "18.14"
$ws.Range("E27").Value = "  -3.78%  "
$ws.Range("D28").Value = "5.067"
$ws.Range("E28").Value = "  -6.70%  "
$ws.Range("D29").Value = "113.44"
$ws.Range("E29").Value = "  -3.42%  "
$ws.Range("D30").Value = "1.689"
$ws.Range("E30").Value = "  -9.73%  "
$ws.Range("D31").Value = "0.08897"
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("D32").Value = "0.7397"
$ws.Range("E32").Value = "  -7.19%  "
$ws.Range("D33").Value = "2.936"
$ws.Range("E33").Value = "  -0.34%  "
$ws.Range("D34").Value = "4.441"
$ws.Range("E34").Value = "  -5.87%  "
$ws.Range("D35").Value = "1.107"
$ws.Range("E35").Value = "  -6.97%  "
$ws.Range("E36").Value = "  +0.04%  "
$ws.Range("D37").Value = "1.070"
$ws.Range("E37").Value = "  -5.28%  "
$ws.Range("D38").Value = "0.05192"
$ws.Range("E38").Value = "  -5.09%  "
$ws.Range("D39").Value = "0.01905"
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("D40").Value = "0.1634"
$ws.Range("E40").Value = "  -3.33%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.705"
$ws.Range("E41").Value = "  -6.28%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.4968"
$ws.Range("E42").Value = "  -4.12%  "
$ws.Range("D43").Value = "6.273"
$ws.Range("E43").Value = "  -9.06%  "
$ws.Range("D44").Value = "8.186"
$ws.Range("E44").Value = "  -5.57%  "
$ws.Range("D45").Value = "105.09"
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "10.21"
$ws.Range("E46").Value = "  -3.78%  "
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  -3.51%  "
$ws.Range("D49").Value = "0.4574"
$ws.Range("E49").Value = "  -4.37%  "
$ws.Range("D50").Value = "1.592"
$ws.Range("E50").Value = "  -4.05%  "
$ws.Range("D51").Value = "62.65"
$ws.Range("E51").Value = "  -4.73%  "
